$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.539.56'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.318.68'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  +0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '515.26'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.64%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '131.60'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -1.75%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +0.36%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.532'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -0.75%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.100'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -2.23%  '
$ws.Range('E10').Value = '  +0.18%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '5.23'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.64%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.337'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -1.23%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.739.04'
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '23.47'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').Value = '56.560.75'
$ws.Range('E15').Value = '  +0.00%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.0000132'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('D17').Value = '2.333.41'
$ws.Range('E17').Value = '  +0.42%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '10.41'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.20%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '330.50'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +2.52%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.15'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.52%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.71'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +2.46%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.996'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.37%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '61.10'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.58%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.165'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +1.45%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '8.64'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +9.20%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.01'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +1.25%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '1.31'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +2.09%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '168.47'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +0.92%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.69'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('D30').Value = '0.0₃0719'
$ws.Range('E30').Value = '  -2.63%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '6.15'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.31%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '18.34'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('E33').Value = '  -0.03%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.54%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.24'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -0.19%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '3.93'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -2.10%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.882'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -3.88%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.58'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +1.94%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '38.73'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +2.36%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '148.20'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +6.64%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '288.91'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +3.83%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.374'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -2.12%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '3.60'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +0.39%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '5.09'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -2.05%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0928'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.46%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0496'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -1.68%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.555'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.35%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '18.08'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.70%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.0215'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -0.73%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '17.17'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +2.13%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '11.02'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.90%  '
